$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear previous style/content on rows 5-10 and row 18 (old data no longer present)
$ws.Range("A5:H10").Clear()
$ws.Range("A18:H18").Clear()

# Rewrite rows 1-4: strip style, set new values
$ws.Range("A1:D4").ClearFormats()

$ws.Range("A1").Value = 10
$ws.Range("A2").Value = 20
$ws.Range("A3").Value = 30
$ws.Range("A4").Value = 40

$ws.Range("B1").Value = "item1"
$ws.Range("B2").Value = "item2"
$ws.Range("B3").Value = "item3"
$ws.Range("B4").Value = "item4"

$ws.Range("C1").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1

$ws.Range("D1").Value = "javi"
$ws.Range("D2").Value = "pedro"
$ws.Range("D3").Value = "cris"
$ws.Range("D4").Value = "fran"

# Update selection to match target (D5)
$ws.Range("D5").Select()
